$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Delete the "License Information" Heading2 paragraph entirely.
# ------------------------------------------------------------------
$licPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "License Information") {
        $licPara = $p
        break
    }
}
if ($licPara -eq $null) {
    throw "Could not find 'License Information' paragraph"
}
$licPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Locate the license paragraph (starts with "Translation Questions
#    (unfoldingWord) is based on") and the following paragraph
#    ("This PDF version is provided under the same license."), then
#    merge them into a single paragraph by removing the paragraph
#    mark between them.
# ------------------------------------------------------------------
$basedOnPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.IndexOf("is based on") -ge 0) {
        $basedOnPara = $p
        break
    }
}
if ($basedOnPara -eq $null) {
    throw "Could not find the license paragraph"
}

$paraStart = $basedOnPara.Range.Start

# find end of "This PDF version is provided under the same license."
$pdfRange = $d.Range($basedOnPara.Range.Start, $d.Content.End)
$found = $pdfRange.Find.Execute("This PDF version is provided under the same license.")
if (-not $found) {
    throw "Could not find 'This PDF version...' text"
}
$pdfTextEnd = $pdfRange.End

# Delete the paragraph mark that currently separates the two paragraphs
# (the mark that falls between basedOnPara's end and the "This PDF..."
# paragraph), so the two paragraphs become one. The mark is the
# character right before the "This PDF" paragraph's content range
# starts growing into the previous paragraph; easiest is to just
# delete every paragraph mark between paraStart and pdfTextEnd except
# none - there is exactly one paragraph boundary in between, locate it:
# it's basedOnPara's own trailing paragraph mark.
$boundary = $basedOnPara.Range.End - 1
$markRange = $d.Range($boundary, $boundary + 1)
$markRange.Delete()
# After deleting the mark, the two paragraphs are merged; pdfTextEnd
# shifts left by 1 character.
$pdfTextEnd = $pdfTextEnd - 1

# ------------------------------------------------------------------
# 3) Replace all text content of the merged paragraph, from its start
#    up to (and including) "...same license." with the new runs,
#    leaving the trailing (empty) runs / paragraph mark untouched.
# ------------------------------------------------------------------
$delRange = $d.Range($paraStart, $pdfTextEnd)
$delRange.Delete()

$global:pos = $paraStart

function Insert-Plain([string]$text) {
    $r = $d.Range($global:pos, $global:pos)
    $r.InsertAfter($text)
    $newEnd = $global:pos + $text.Length
    $fr = $d.Range($global:pos, $newEnd)
    $fr.Font.Bold = 0
    $global:pos = $newEnd
}

function Insert-Bold([string]$text) {
    $r = $d.Range($global:pos, $global:pos)
    $r.InsertAfter($text)
    $newEnd = $global:pos + $text.Length
    $fr = $d.Range($global:pos, $newEnd)
    $fr.Font.Bold = 1
    $global:pos = $newEnd
}

Insert-Bold "unfoldingWord® Translation Questions"
Insert-Plain " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. "
Insert-Plain "unfoldingWord® Translation Questions"
Insert-Plain " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from "
Insert-Plain "unfoldingWord® Translation Questions"
Insert-Plain " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

Write-Host "Done. Merged paragraph text:"
Write-Host $d.Range($paraStart, $global:pos).Text
